# Adds a new forecast "vintage" column BB (9th-archive equivalent) to the
# naive-forecaster QoQ GDP table, and one additional forecast-horizon row
# (row 83) to extend the diagonal staircase pattern by one more quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cell BB1 = new vintage date (2025-11-25), same style as
#    the rest of row 1 (B1:BA1) -- bold / centered / bordered date format.
# ---------------------------------------------------------------------
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 54).Value = 45986

# ---------------------------------------------------------------------
# 2. Rows 2-71: new column BB simply repeats column BA's forecast value
#    (the series has already converged/flattened by the last vintage),
#    so copy BA2:BA71 straight across into BB2:BB71.
# ---------------------------------------------------------------------
$ws.Range("BA2:BA71").Copy() | Out-Null
$ws.Range("BB2:BB71").PasteSpecial(-4163) | Out-Null

# ---------------------------------------------------------------------
# 3. Row 72: new vintage's own forecast value differs from BA's.
# ---------------------------------------------------------------------
$ws.Cells.Item(72, 54).Value = -0.2099029780610664

# ---------------------------------------------------------------------
# 4. Rows 73-82: beyond the new vintage's forecast horizon -> 0.
# ---------------------------------------------------------------------
$ws.Cells.Item(73, 54).Value = 0
$ws.Cells.Item(74, 54).Value = 0
$ws.Cells.Item(75, 54).Value = 0
$ws.Cells.Item(76, 54).Value = 0
$ws.Cells.Item(77, 54).Value = 0
$ws.Cells.Item(78, 54).Value = 0
$ws.Cells.Item(79, 54).Value = 0
$ws.Cells.Item(80, 54).Value = 0
$ws.Cells.Item(81, 54).Value = 0
$ws.Cells.Item(82, 54).Value = 0

# ---------------------------------------------------------------------
# 5. New row 83: one more quarter-end target period (2028-06-30), with
#    the same style as the rest of column A, plus its BB entry (0, since
#    it's beyond the new vintage's forecast horizon too).
# ---------------------------------------------------------------------
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(83, 1).Value = 46934
$ws.Cells.Item(83, 54).Value = 0
